$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1549815498154982
$ws.Range("C2").Value = 0.6199261992619927
$ws.Range("J2").Value = 0.02952029520295203
$ws.Range("P2").Value = 0.1107011070110701
$ws.Range("S2").Value = 0.08487084870848709
$ws.Range("C3").Value = 0.01685393258426966
$ws.Range("J3").Value = 0.08426966292134831
$ws.Range("P3").Value = 0.6853932584269663
$ws.Range("S3").Value = 0.2134831460674157
$ws.Range("J4").Value = 0.05172413793103448
$ws.Range("P4").Value = 0.7068965517241379
$ws.Range("S4").Value = 0.2413793103448276
$ws.Range("B6").Value = 0.03846153846153846
$ws.Range("D6").Value = 0.01923076923076923
$ws.Range("F6").Value = 0.04326923076923077
$ws.Range("J6").Value = 0.3173076923076923
$ws.Range("O6").Value = 0.02403846153846154
$ws.Range("Q6").Value = 0.1682692307692308
$ws.Range("R6").Value = 0.07211538461538461
$ws.Range("S6").Value = 0.3173076923076923
$ws.Range("B7").Value = 0.06422018348623854
$ws.Range("D7").Value = 0.01834862385321101
$ws.Range("E7").Value = 0.009174311926605505
$ws.Range("F7").Value = 0.06880733944954129
$ws.Range("J7").Value = 0.2018348623853211
$ws.Range("O7").Value = 0.02293577981651376
$ws.Range("Q7").Value = 0.1467889908256881
$ws.Range("R7").Value = 0.1100917431192661
$ws.Range("S7").Value = 0.3577981651376147
$ws.Range("B8").Value = 0.1023890784982935
$ws.Range("D8").Value = 0.02389078498293516
$ws.Range("F8").Value = 0.05119453924914676
$ws.Range("J8").Value = 0.1467576791808874
$ws.Range("O8").Value = 0.0273037542662116
$ws.Range("Q8").Value = 0.1945392491467577
$ws.Range("R8").Value = 0.09215017064846416
$ws.Range("S8").Value = 0.3617747440273038
$ws.Range("B9").Value = 0.1161825726141079
$ws.Range("D9").Value = 0.02074688796680498
$ws.Range("F9").Value = 0.05394190871369295
$ws.Range("J9").Value = 0.1286307053941909
$ws.Range("O9").Value = 0.04149377593360996
$ws.Range("Q9").Value = 0.1825726141078838
$ws.Range("R9").Value = 0.08713692946058091
$ws.Range("S9").Value = 0.3692946058091287
$ws.Range("B10").Value = 0.1072243346007605
$ws.Range("D10").Value = 0.02965779467680608
$ws.Range("F10").Value = 0.05779467680608365
$ws.Range("J10").Value = 0.1422053231939163
$ws.Range("O10").Value = 0.02129277566539924
$ws.Range("Q10").Value = 0.2387832699619772
$ws.Range("R10").Value = 0.07680608365019011
$ws.Range("S10").Value = 0.3262357414448669
$ws.Range("G11").Value = 0.1238095238095238
$ws.Range("J11").Value = 0.08888888888888889
$ws.Range("K11").Value = 0.1714285714285714
$ws.Range("L11").Value = 0.6158730158730159
$ws.Range("G12").Value = 0.775
$ws.Range("J12").Value = 0.185
$ws.Range("K12").Value = 0.005
$ws.Range("L12").Value = 0.02
$ws.Range("S12").Value = 0.015
$ws.Range("F13").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.7368421052631579
$ws.Range("J13").Value = 0.2368421052631579
$ws.Range("F15").Value = 0.01477832512315271
$ws.Range("H15").Value = 0.1330049261083744
$ws.Range("I15").Value = 0.1182266009852217
$ws.Range("J15").Value = 0.3842364532019704
$ws.Range("K15").Value = 0.03940886699507389
$ws.Range("M15").Value = 0.004926108374384237
$ws.Range("O15").Value = 0.06896551724137931
$ws.Range("S15").Value = 0.2364532019704433
$ws.Range("F16").Value = 0.02645502645502645
$ws.Range("H16").Value = 0.1216931216931217
$ws.Range("I16").Value = 0.1164021164021164
$ws.Range("J16").Value = 0.4814814814814815
$ws.Range("K16").Value = 0.1164021164021164
$ws.Range("M16").Value = 0.005291005291005291
$ws.Range("O16").Value = 0.03703703703703703
$ws.Range("S16").Value = 0.09523809523809523
$ws.Range("F17").Value = 0.02505219206680585
$ws.Range("H17").Value = 0.1461377870563674
$ws.Range("I17").Value = 0.1148225469728601
$ws.Range("J17").Value = 0.4154488517745303
$ws.Range("K17").Value = 0.1377870563674322
$ws.Range("M17").Value = 0.01670146137787056
$ws.Range("N17").Value = 0.006263048016701462
$ws.Range("O17").Value = 0.04801670146137787
$ws.Range("S17").Value = 0.08977035490605428
$ws.Range("F18").Value = 0.0213903743315508
$ws.Range("H18").Value = 0.1122994652406417
$ws.Range("I18").Value = 0.09090909090909091
$ws.Range("J18").Value = 0.5187165775401069
$ws.Range("K18").Value = 0.1016042780748663
$ws.Range("M18").Value = 0.0374331550802139
$ws.Range("N18").Value = 0.0053475935828877
$ws.Range("O18").Value = 0.0374331550802139
$ws.Range("S18").Value = 0.0748663101604278
$ws.Range("F19").Value = 0.02783109404990403
$ws.Range("H19").Value = 0.1497120921305182
$ws.Range("I19").Value = 0.1151631477927063
$ws.Range("J19").Value = 0.3973128598848368
$ws.Range("K19").Value = 0.1314779270633397
$ws.Range("M19").Value = 0.01919385796545105
$ws.Range("N19").Value = 0.002879078694817658
$ws.Range("O19").Value = 0.07005758157389635
$ws.Range("S19").Value = 0.08637236084452975
